$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing "res" values in column C ---
# Row 3: Invalid -> Valid
$ws.Range("C3").Value = "Valid"
# Row 4: Invalid -> Valid
$ws.Range("C4").Value = "Valid"

# Row 6: C6 loses its green/fill style, now matches the plain "bold, no fill" style (same as C2)
$ws.Range("C6").Style = $ws.Range("C2").Style

# --- Add a new row 7 ---
# A7: empty cell, but carries the same style as C2/B3 (plain bold, no fill)
$ws.Range("A7").Style = $ws.Range("C2").Style

# B7: "test@123" with a mailto hyperlink, same style as A2 (hyperlink, no fill)
$ws.Range("B7").Style = $ws.Range("A2").Style
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:test@123", "", "", "test@123")

# C7: "Invalid", same style as C2 (plain bold, no fill)
$ws.Range("C7").Style = $ws.Range("C2").Style
$ws.Range("C7").Value = "Invalid"

# Update the selection to match the final cursor position
$ws.Range("C7").Select()
